$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 fresh rows right above row 596, pushing the existing
# rows 596-604 down to 600-608 (values/formats move with them).
$ws.Range("A596:R599").EntireRow.Insert()

# Row 596 - new weekly record (Extra)
$ws.Cells.Item(596, 1).Value = 4
$ws.Cells.Item(596, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(596, 3).Value = "Los Lagos"
$ws.Cells.Item(596, 4).Value = 44656
$ws.Cells.Item(596, 5).Value = 10
$ws.Cells.Item(596, 6).Value = 100112020
$ws.Cells.Item(596, 7).Value = "Tomate"
$ws.Cells.Item(596, 8).Value = "Larga vida"
$ws.Cells.Item(596, 9).Value = "Extra"
$ws.Cells.Item(596, 10).Value = 300
$ws.Cells.Item(596, 11).Value = 25000
$ws.Cells.Item(596, 12).Value = 25000
$ws.Cells.Item(596, 13).Value = 25000
$ws.Cells.Item(596, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(596, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(596, 16).Value = 1389
$ws.Cells.Item(596, 17).Value = 18
$ws.Cells.Item(596, 18).Value = "Hortaliza"

# Row 597 - new weekly record (Extra)
$ws.Cells.Item(597, 1).Value = 4
$ws.Cells.Item(597, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(597, 3).Value = "Los Lagos"
$ws.Cells.Item(597, 4).Value = 44656
$ws.Cells.Item(597, 5).Value = 10
$ws.Cells.Item(597, 6).Value = 100112020
$ws.Cells.Item(597, 7).Value = "Tomate"
$ws.Cells.Item(597, 8).Value = "Larga vida"
$ws.Cells.Item(597, 9).Value = "Extra"
$ws.Cells.Item(597, 10).Value = 500
$ws.Cells.Item(597, 11).Value = 26000
$ws.Cells.Item(597, 12).Value = 26500
$ws.Cells.Item(597, 13).Value = 26250
$ws.Cells.Item(597, 14).Value = "`$/bandeja 20 kilos"
$ws.Cells.Item(597, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(597, 16).Value = 1312
$ws.Cells.Item(597, 17).Value = 20
$ws.Cells.Item(597, 18).Value = "Hortaliza"

# Row 598 - new weekly record (Primera)
$ws.Cells.Item(598, 1).Value = 4
$ws.Cells.Item(598, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(598, 3).Value = "Los Lagos"
$ws.Cells.Item(598, 4).Value = 44656
$ws.Cells.Item(598, 5).Value = 10
$ws.Cells.Item(598, 6).Value = 100112020
$ws.Cells.Item(598, 7).Value = "Tomate"
$ws.Cells.Item(598, 8).Value = "Larga vida"
$ws.Cells.Item(598, 9).Value = "Primera"
$ws.Cells.Item(598, 10).Value = 300
$ws.Cells.Item(598, 11).Value = 23000
$ws.Cells.Item(598, 12).Value = 23000
$ws.Cells.Item(598, 13).Value = 23000
$ws.Cells.Item(598, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(598, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(598, 16).Value = 1278
$ws.Cells.Item(598, 17).Value = 18
$ws.Cells.Item(598, 18).Value = "Hortaliza"

# Row 599 - new weekly record (Tercera)
$ws.Cells.Item(599, 1).Value = 4
$ws.Cells.Item(599, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(599, 3).Value = "Los Lagos"
$ws.Cells.Item(599, 4).Value = 44656
$ws.Cells.Item(599, 5).Value = 10
$ws.Cells.Item(599, 6).Value = 100112020
$ws.Cells.Item(599, 7).Value = "Tomate"
$ws.Cells.Item(599, 8).Value = "Larga vida"
$ws.Cells.Item(599, 9).Value = "Tercera"
$ws.Cells.Item(599, 10).Value = 300
$ws.Cells.Item(599, 11).Value = 17000
$ws.Cells.Item(599, 12).Value = 17000
$ws.Cells.Item(599, 13).Value = 17000
$ws.Cells.Item(599, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(599, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(599, 16).Value = 944
$ws.Cells.Item(599, 17).Value = 18
$ws.Cells.Item(599, 18).Value = "Hortaliza"
